$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right after "Source:" (row 30), before "DUBAI STATISTICS CENTER"
# (old row 31). This shifts old rows 31-37 down to 32-38.
$ws.Rows("31:31").Insert()

# After the insert:
#   row 31 = blank (new)
#   row 32 = "DUBAI STATISTICS CENTER" (was row 31)
#   row 33 = URL, still hyperlink-styled (was row 32)
#   row 34 = blank (was row 33)
#   row 37 = "DSC" (was row 36)
#   row 38 = old citation text (was row 37)
#
# The target layout wants the blank row and the URL row swapped (blank at 33,
# URL at 34), with the URL cell no longer hyperlinked, so fix that up here.

$ws.Range("A31").Style = "source"
$ws.Range("A31").Value = ""

$ws.Range("A32").Style = "source"

# Remove the hyperlink that was attached to the URL cell before the row insert.
$ws.Range("A33").Hyperlinks.Delete()

$url = "http://dsc.gov.ae/Reports/Establishments%20Distribution%20%20by%20Economic%20Activity%20and%20Workers%20Group.pdf"

# Row 33 becomes blank, row 34 gets the URL text (plain, no hyperlink).
$ws.Range("A33").Value = ""
$ws.Range("A33").Style = "source"

$ws.Range("A34").Value = $url
$ws.Range("A34").Style = "source"

# Row 37 keeps "DSC" already in place (unchanged by the insert, just moved).
$ws.Range("A37").Style = "title"

# Row 38 gets the new, longer citation text (replaces the old citation text).
$newCitation = "Dubai SME an agency of the Department of Economic Development, ""THE STATE OF SMALL & MEDIUM ENTERPRISES (SMEs) IN DUBAI"". Available at http://www.sme.ae/upload/category/SME_Report_English.pdf. `nMohammed Bin Rashid - Establishment for SME Development, ""THE DEFINITION OF MICRO, SMALL & MEDIUM ENTERPRISES (MSMEs) OF DUBAI"". Available at http://www.sme.ae/upload/category/SME_l_DEFINITION_l_English.pdf`nIt is important to note that UAE Issued Federal Law No. 2 of 2014 for Small and Medium Enterprises. In it, it is detailed in article 4 that the UAE Cabinet ""must issue a decree on the united definition of SMEs, based on the recommedations of the SME Council, in consultation with competent local authorities""."

$ws.Range("A38").Value = $newCitation
$ws.Range("A38").Style = "source"
